$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Sema4a"
$row2[0,2] = "Plxnb1"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 25.81602566666667
$row2[0,7] = 77.44807700000001
$row2[0,8] = 0.7742517153725241
$row2[0,9] = 0.7742517153725241
$row2[0,10] = 2
$row2[0,11] = 0.6666666666666666
$row2[0,12] = 0.6893159999999999
$row2[0,13] = 2.067948
$row2[0,14] = 0.1083604551316437
$row2[0,15] = 0.1083604551316437
$row2[0,16] = 17.795399548444
$row2[0,17] = 160.158595935996
$row2[0,18] = 0.08389826826422259
$row2[0,19] = 0.08389826826422257
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Sema4a"
$row3[0,2] = "Plxnb1"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 25.81602566666667
$row3[0,7] = 77.44807700000001
$row3[0,8] = 0.7742517153725241
$row3[0,9] = 0.7742517153725241
$row3[0,10] = 2
$row3[0,11] = 0.6666666666666666
$row3[0,12] = 0.082959
$row3[0,13] = 0.248877
$row3[0,14] = 0.01304115238477858
$row3[0,15] = 0.01304115238477858
$row3[0,16] = 2.141671673281
$row3[0,17] = 19.275045059529
$row3[0,18] = 0.0100971346043493
$row3[0,19] = 0.0100971346043493
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Sema4a"
$row4[0,2] = "Plxnb1"
$row4[0,3] = "sCs"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 25.81602566666667
$row4[0,7] = 77.44807700000001
$row4[0,8] = 0.7742517153725241
$row4[0,9] = 0.7742517153725241
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 5.589049333333333
$row4[0,13] = 16.767148
$row4[0,14] = 0.8785983924835777
$row4[0,15] = 0.8785983924835776
$row4[0,16] = 144.2870410415996
$row4[0,17] = 1298.583369374396
$row4[0,18] = 0.6802563125039522
$row4[0,19] = 0.6802563125039522
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "FAPs"
$row5[0,1] = "Sema4a"
$row5[0,2] = "Plxnb1"
$row5[0,3] = "ECs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 2.718527666666667
$row5[0,7] = 8.155583
$row5[0,8] = 0.08153170965901445
$row5[0,9] = 0.08153170965901445
$row5[0,10] = 2
$row5[0,11] = 0.6666666666666666
$row5[0,12] = 0.6893159999999999
$row5[0,13] = 2.067948
$row5[0,14] = 0.1083604551316437
$row5[0,15] = 0.1083604551316437
$row5[0,16] = 1.873924617076
$row5[0,17] = 16.865321553684
$row5[0,18] = 0.008834813166311838
$row5[0,19] = 0.008834813166311836
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Sema4a"
$row6[0,2] = "Plxnb1"
$row6[0,3] = "FAPs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 2.718527666666667
$row6[0,7] = 8.155583
$row6[0,8] = 0.08153170965901445
$row6[0,9] = 0.08153170965901445
$row6[0,10] = 2
$row6[0,11] = 0.6666666666666666
$row6[0,12] = 0.082959
$row6[0,13] = 0.248877
$row6[0,14] = 0.01304115238477858
$row6[0,15] = 0.01304115238477858
$row6[0,16] = 0.225526336699
$row6[0,17] = 2.029737030291
$row6[0,18] = 0.001063267449854731
$row6[0,19] = 0.001063267449854731
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Sema4a"
$row7[0,2] = "Plxnb1"
$row7[0,3] = "sCs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 2.718527666666667
$row7[0,7] = 8.155583
$row7[0,8] = 0.08153170965901445
$row7[0,9] = 0.08153170965901445
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 5.589049333333333
$row7[0,13] = 16.767148
$row7[0,14] = 0.8785983924835777
$row7[0,15] = 0.8785983924835776
$row7[0,16] = 15.19398524303155
$row7[0,17] = 136.745867187284
$row7[0,18] = 0.07163362904284788
$row7[0,19] = 0.07163362904284787
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "sCs"
$row8[0,1] = "Sema4a"
$row8[0,2] = "Plxnb1"
$row8[0,3] = "ECs"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 4.808641333333334
$row8[0,7] = 14.425924
$row8[0,8] = 0.1442165749684613
$row8[0,9] = 0.1442165749684613
$row8[0,10] = 2
$row8[0,11] = 0.6666666666666666
$row8[0,12] = 0.6893159999999999
$row8[0,13] = 2.067948
$row8[0,14] = 0.1083604551316437
$row8[0,15] = 0.1083604551316437
$row8[0,16] = 3.314673409328
$row8[0,17] = 29.832060683952
$row8[0,18] = 0.01562737370110928
$row8[0,19] = 0.01562737370110928
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "sCs"
$row9[0,1] = "Sema4a"
$row9[0,2] = "Plxnb1"
$row9[0,3] = "FAPs"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 4.808641333333334
$row9[0,7] = 14.425924
$row9[0,8] = 0.1442165749684613
$row9[0,9] = 0.1442165749684613
$row9[0,10] = 2
$row9[0,11] = 0.6666666666666666
$row9[0,12] = 0.082959
$row9[0,13] = 0.248877
$row9[0,14] = 0.01304115238477858
$row9[0,15] = 0.01304115238477858
$row9[0,16] = 0.3989200763720001
$row9[0,17] = 3.590280687348001
$row9[0,18] = 0.001880750330574548
$row9[0,19] = 0.001880750330574548
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "sCs"
$row10[0,1] = "Sema4a"
$row10[0,2] = "Plxnb1"
$row10[0,3] = "sCs"
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 4.808641333333334
$row10[0,7] = 14.425924
$row10[0,8] = 0.1442165749684613
$row10[0,9] = 0.1442165749684613
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 5.589049333333333
$row10[0,13] = 16.767148
$row10[0,14] = 0.8785983924835777
$row10[0,15] = 0.8785983924835776
$row10[0,16] = 26.87573363830578
$row10[0,17] = 241.881602744752
$row10[0,18] = 0.1267084509367775
$row10[0,19] = 0.1267084509367775
$ws.Range("A10:T10").Value = $row10
